$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename mapping values: "allow_entry" split into distinct initial/final targets
$ws.Range("C3").Value = "allow_entry_init"
$ws.Range("C4").Value = "allow_entry_final"

# Move the active selection to D2 (matches saved cursor position)
$ws.Range("D2").Select()

# Page setup touched (portrait orientation) as part of the save
$ws.PageSetup.Orientation = 1
